$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit duplicates the value of A5 ("EV004") into a new cell C5,
# extending the sheet's used range from A1:A26 to A1:C26.
$ws.Range("C5").Value = $ws.Range("A5").Value()

# Excel recomputes the per-row "spans" hint in 16-row blocks based on the
# furthest populated column seen anywhere in that block. Touch a cell's
# formatting (no value) in the second block (rows 17-32) so that block's
# rows also report the widened 1:3 span, matching the rest of the sheet -
# without leaving any real value/content in that cell.
$ws.Range("Z1").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0
